$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.373.75'
$ws.Range("E2").Value = '  +0.05%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.672.83'
$ws.Range("E3").Value = '  -0.35%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '645.47'
$ws.Range("E5").Value = '  -5.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.42'
$ws.Range("E6").Value = '  -0.86%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  +0.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.145'
$ws.Range("E9").Value = '  -0.91%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.10'
$ws.Range("E10").Value = '  -0.69%  '
$ws.Range("E11").Value = '  +1.20%  '
$ws.Range("E12").Value = '  -0.69%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.291.86'
$ws.Range("E13").Value = '  -0.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.47'
$ws.Range("E14").Value = '  +0.45%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.687.94'
$ws.Range("E15").Value = '  +0.28%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '69.355.38'
$ws.Range("E16").Value = '  -0.01%  '
$ws.Range("E17").Value = '  -0.09%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '15.83'
$ws.Range("E18").Value = '  -1.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.46'
$ws.Range("E19").Value = '  +0.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '468.23'
$ws.Range("E20").Value = '  -0.78%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.99'
$ws.Range("E21").Value = '  +1.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.644'
$ws.Range("E22").Value = '  -0.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '79.29'
$ws.Range("E23").Value = '  -0.87%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.817.86'
$ws.Range("E24").Value = '  -0.43%  '
$ws.Range("E26").Value = '  +0.93%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.69'
$ws.Range("E27").Value = '  -1.64%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.96'
$ws.Range("E28").Value = '  -1.35%  '
$ws.Range("E29").Value = '  -3.09%  '
$ws.Range("E30").Value = '  -2.32%  '
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.00'
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("B32").Value = 'Binance-PegBSC-USD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("E32").Value = '  -0.25%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.96'
$ws.Range("E33").Value = '  +0.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.39'
$ws.Range("E34").Value = '  -2.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.664.74'
$ws.Range("E35").Value = '  -0.26%  '
$ws.Range("E36").Value = '  +0.20%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.38'
$ws.Range("E37").Value = '  +0.21%  '
$ws.Range("E39").Value = '  +5.41%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.86'
$ws.Range("E40").Value = '  -6.13%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  -0.17%  '
$ws.Range("E42").Value = '  -1.91%  '
$ws.Range("E43").Value = '  -1.64%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.924'
$ws.Range("E44").Value = '  -1.75%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '46.95'
$ws.Range("E45").Value = '  +0.17%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '28.10'
$ws.Range("E46").Value = '  -1.65%  '
$ws.Range("B47").Value = 'dogwifhat'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.69'
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("E48").Value = '  -3.29%  '
$ws.Range("B49").Value = 'FLOKI'
$ws.Range("C49").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000265'
$ws.Range("E49").Value = '  -4.01%  '
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.78'
$ws.Range("E50").Value = '  -0.57%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.23'
$ws.Range("E51").Value = '  -3.90%  '
